# Update "想去人数" (want-to-go count) and "最低票价" (lowest price) figures
# across the "展览" and "全部类型" worksheets, as published at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4510
$ws1.Range("G5").Value = 55
$ws1.Range("G6").Value = 40
$ws1.Range("F7").Value = 154
$ws1.Range("F8").Value = 635
$ws1.Range("F11").Value = 1360
$ws1.Range("F12").Value = 26
$ws1.Range("F13").Value = 2968
$ws1.Range("F15").Value = 666

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4510
$ws4.Range("G5").Value = 55
$ws4.Range("G6").Value = 40
$ws4.Range("F7").Value = 154
$ws4.Range("F8").Value = 635
$ws4.Range("F12").Value = 1360
$ws4.Range("F13").Value = 26
$ws4.Range("F14").Value = 2968
$ws4.Range("F16").Value = 666

$wb.Save()
